$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns D and E sized to fit their new contents (matches the author's manual resize)
$ws.Columns.Item(4).ColumnWidth = 54.736979166666664
$ws.Columns.Item(5).ColumnWidth = 79.73697916666667

# Fill in the new Q&A notes, in the order the shared strings were authored
$ws.Range("D8").Value = "adnotacje sa odczytywane przez srodowiska, frameworki. "
$ws.Range("D8").Font.Name = "Arial"

$ws.Range("D9").Value = "iterface maja tylko przepisy na metody"
$ws.Range("D9").Font.Name = "Arial"

$ws.Range("C19").Value = "implementacja"

$ws.Range("D13").Value = "late binding"
$ws.Range("D13").Font.Name = "Arial"

$ws.Range("E13").Value = "switch random. Tablica trzech ksztaltow. Losuje i wrzycasz trojkat, okrag, kwadrat. "
$ws.Range("E13").Font.Name = "Arial"

$ws.Range("E12").Value = "overide"
$ws.Range("E12").Font.Name = "Arial"

# Move the active selection the way the author left it
$ws.Range("D18").Select() | Out-Null
